# Insert a new weekly price-report row at row 530 (pushing the existing
# rows 530-630 down to 531-631, dimension grows from A1:R630 to A1:R631),
# then populate the new row with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 530; formatting/styles of the
# row are inherited the same way Excel's own Insert does.
$ws.Rows.Item(530).Insert()

$row = 530

$ws.Cells.Item($row, 1).Value  = 5
$ws.Cells.Item($row, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item($row, 3).Value  = "Maule"
$ws.Cells.Item($row, 4).Value  = 45015
$ws.Cells.Item($row, 5).Value  = 7
$ws.Cells.Item($row, 6).Value  = 100112043
$ws.Cells.Item($row, 7).Value  = "Pepino ensalada"
$ws.Cells.Item($row, 8).Value  = "Sin especificar"
$ws.Cells.Item($row, 9).Value  = "Primera"
$ws.Cells.Item($row, 10).Value = 300
$ws.Cells.Item($row, 11).Value = 6000
$ws.Cells.Item($row, 12).Value = 6000
$ws.Cells.Item($row, 13).Value = 6000
$ws.Cells.Item($row, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 100
$ws.Cells.Item($row, 17).Value = 60
$ws.Cells.Item($row, 18).Value = "Hortaliza"
